$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from A10 (the previous year-label cell) onto the new A11 label cell
# so the new cell reuses the existing bold/centered/bordered style instead of
# creating a brand-new style entry.
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)

# Row 11: 2021 data
$ws.Range("A11").Value = "2021年"
$ws.Range("B11").Value = 28510.26
$ws.Range("C11").Value = 7191.92
$ws.Range("D11").Value = 2547.91
$ws.Range("F11").Value = 14055.93
$ws.Range("G11").Value = 56939.13
$ws.Range("H11").Value = 7987.7
$ws.Range("I11").Value = 32899.32
$ws.Range("J11").Value = 2916.38
$ws.Range("K11").Value = 995876.22
$ws.Range("L11").Value = 3898.53
$ws.Range("M11").Value = 1344.62
$ws.Range("N11").Value = 2849.52
$ws.Range("O11").Value = 4982.68
$ws.Range("P11").Value = 34564.18
$ws.Range("Q11").Value = 4366.88
$ws.Range("R11").Value = 1411.68
$ws.Range("S11").Value = 13213.3
$ws.Range("T11").Value = 12456.73
$ws.Range("U11").Value = 70764.1
$ws.Range("V11").Value = 11484.35
$ws.Range("W11").Value = 59253.94
$ws.Range("X11").Value = 6832.09
$ws.Range("Y11").Value = 125233.84
$ws.Range("Z11").Value = 63675.07
$ws.Range("AA11").Value = 3569.39
$ws.Range("AB11").Value = 37544.53
$ws.Range("AC11").Value = 20785.56
$ws.Range("AD11").Value = 11046.52
$ws.Range("AE11").Value = 6671.61
$ws.Range("AF11").Value = 127855.9
$ws.Range("AG11").Value = 33385.49
$ws.Range("AH11").Value = 11001.54
$ws.Range("AI11").Value = 14285.46
$ws.Range("AJ11").Value = 1941.02
$ws.Range("AK11").Value = 18503.76
$ws.Range("AL11").Value = 23170.62
$ws.Range("AM11").Value = 30330.17
$ws.Range("AN11").Value = 1794.27
$ws.Range("AO11").Value = 12002.36
$ws.Range("AP11").Value = 62740.28
$ws.Range("AQ11").Value = 9860.74
